$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 28.545454
$ws.Range("I11").Value = 28.545454
$ws.Range("K11").Value = 28.545454
$ws.Range("M11").Value = 111.454546

$ws.Range("H17").Value = 3628520
$ws.Range("J17").Value = 3786007.8
$ws.Range("L17").Value = 11358023.4
$ws.Range("N17").Value = -11358359.4

$ws.Range("H28").Value = 183.58824
$ws.Range("I28").Value = 183.58824
$ws.Range("K28").Value = 183.58824
$ws.Range("M28").Value = 301.41176

$ws.Range("H47").Value = 25014
$ws.Range("I47").Value = 25028.5
$ws.Range("K47").Value = 25028.5
$ws.Range("M47").Value = -24056.5

$ws.Range("H51").Value = 6618.7915
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 6618.7915
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6618.7915
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7586.7915

$ws.Range("H86").Value = 3542.7273
$ws.Range("I86").Value = 3432.1667
$ws.Range("J86").Value = 3675.4
$ws.Range("K86").Value = 3432.1667
$ws.Range("L86").Value = 3675.4
$ws.Range("M86").Value = -2309.1667
$ws.Range("N86").Value = -5921.4

$ws.Range("H89").Value = 3542.7273
$ws.Range("I89").Value = 3432.1667
$ws.Range("J89").Value = 3675.4
$ws.Range("K89").Value = 17160.8335
$ws.Range("L89").Value = 18377
$ws.Range("M89").Value = -11544.8335
$ws.Range("N89").Value = -29609

$ws.Range("H98").Value = 1736.7916
$ws.Range("I98").Value = 1475.1428
$ws.Range("K98").Value = 1475.1428
$ws.Range("M98").Value = 22.85719999999992

$ws.Range("H111").Value = 1186.7273
$ws.Range("I111").Value = 1186.7273
$ws.Range("K111").Value = 3560.1819
$ws.Range("M111").Value = -493.1819

$ws.Range("H113").Value = 4350.4287
$ws.Range("I113").Value = 3355.0908
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 3355.0908
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -101.0907999999999
$ws.Range("N113").Value = -14508

$ws.Range("H116").Value = 7959.4614
$ws.Range("J116").Value = 9335
$ws.Range("L116").Value = 9335
$ws.Range("N116").Value = -16219

$ws.Range("H121").Value = 5000
$ws.Range("J121").Value = 5000
$ws.Range("L121").Value = 15000
$ws.Range("N121").Value = -18494

$ws.Range("H122").Value = 1736.7916
$ws.Range("I122").Value = 1475.1428
$ws.Range("K122").Value = 4425.428400000001
$ws.Range("M122").Value = -1975.428400000001

$ws.Range("H127").Value = 2332.3333
$ws.Range("I127").Value = 2332.3333
$ws.Range("K127").Value = 6996.999899999999
$ws.Range("M127").Value = -2036.999899999999

$ws.Range("H137").Value = 10351.106
$ws.Range("I137").Value = 4663.524
$ws.Range("K137").Value = 13990.572
$ws.Range("M137").Value = -11440.572

$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

$ws.Range("H141").Value = 3837.125
$ws.Range("I141").Value = 4219.4
$ws.Range("J141").Value = 3200
$ws.Range("K141").Value = 12658.2
$ws.Range("L141").Value = 9600
$ws.Range("M141").Value = -7478.199999999999
$ws.Range("N141").Value = -19960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3029.8823
$ws.Range("I2").Value = 2876.0588
$ws.Range("K2").Value = 2876.0588
$ws.Range("M2").Value = -2763.0588

$ws.Range("H37").Value = 31812.637
$ws.Range("J37").Value = 31812.637
$ws.Range("L37").Value = 31812.637
$ws.Range("N37").Value = -32358.637

$ws.Range("H61").Value = 9175.806
$ws.Range("I61").Value = 7634.1904
$ws.Range("K61").Value = 7634.1904
$ws.Range("M61").Value = -7422.1904

$ws.Range("H63").Value = 4213.5713
$ws.Range("I63").Value = 2297.8
$ws.Range("J63").Value = 9003
$ws.Range("K63").Value = 2297.8
$ws.Range("L63").Value = 9003
$ws.Range("M63").Value = -1611.8
$ws.Range("N63").Value = -10375

$ws.Range("H66").Value = 4213.5713
$ws.Range("I66").Value = 2297.8
$ws.Range("J66").Value = 9003
$ws.Range("K66").Value = 11489
$ws.Range("L66").Value = 45015
$ws.Range("M66").Value = -8057
$ws.Range("N66").Value = -51879

$ws.Range("H74").Value = 18579.5
$ws.Range("I74").Value = 18904.908
$ws.Range("K74").Value = 18904.908
$ws.Range("M74").Value = -18030.908

$ws.Range("H77").Value = 18579.5
$ws.Range("I77").Value = 18904.908
$ws.Range("K77").Value = 94524.53999999999
$ws.Range("M77").Value = -90156.53999999999

$ws.Range("H110").Value = 1304.0834
$ws.Range("I110").Value = 1261.1111
$ws.Range("J110").Value = 1433
$ws.Range("K110").Value = 1261.1111
$ws.Range("L110").Value = 1433
$ws.Range("M110").Value = 783.8888999999999
$ws.Range("N110").Value = -5523

$ws.Range("H116").Value = 3029.8823
$ws.Range("I116").Value = 2876.0588
$ws.Range("K116").Value = 2876.0588
$ws.Range("M116").Value = -582.0587999999998

$ws.Range("H122").Value = 2014.875
$ws.Range("J122").Value = 1889.75
$ws.Range("L122").Value = 5669.25
$ws.Range("N122").Value = -10569.25

$ws.Range("H125").Value = 119999.5
$ws.Range("J125").Value = 119999.5
$ws.Range("L125").Value = 119999.5
$ws.Range("N125").Value = -129839.5

$ws.Range("H132").Value = 2465.147
$ws.Range("I132").Value = 1845.0588
$ws.Range("K132").Value = 5535.1764
$ws.Range("M132").Value = -3005.1764

$ws.Range("H136").Value = 9175.806
$ws.Range("I136").Value = 7634.1904
$ws.Range("K136").Value = 22902.5712
$ws.Range("M136").Value = -20352.5712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3029.8823
$ws.Range("I3").Value = 2876.0588
$ws.Range("K3").Value = 2876.0588
$ws.Range("M3").Value = -2762.0588

$ws.Range("H75").Value = 30235
$ws.Range("I75").Value = 30235
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 30235
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -29299
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 30235
$ws.Range("I78").Value = 30235
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 90705
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -86025
$ws.Range("N78").ClearContents()

$ws.Range("H105").Value = 4406.091
$ws.Range("I105").Value = 3672.7646
$ws.Range("K105").Value = 3672.7646
$ws.Range("M105").Value = -1925.7646

$ws.Range("H107").Value = 1303.8572
$ws.Range("I107").Value = 386.9091
$ws.Range("J107").Value = 4666
$ws.Range("K107").Value = 386.9091
$ws.Range("L107").Value = 4666
$ws.Range("M107").Value = 1533.0909
$ws.Range("N107").Value = -8506

$ws.Range("H134").Value = 8436.4
$ws.Range("I134").Value = 3894.484
$ws.Range("K134").Value = 11683.452
$ws.Range("M134").Value = -9148.451999999999

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5791.222
$ws.Range("J31").Value = 5795.4
$ws.Range("L31").Value = 5795.4
$ws.Range("N31").Value = -6385.4

$ws.Range("H34").Value = 5791.222
$ws.Range("J34").Value = 5795.4
$ws.Range("L34").Value = 5795.4
$ws.Range("N34").Value = -6199.4

$ws.Range("H58").Value = 5071.12
$ws.Range("I58").Value = 3439.6365
$ws.Range("K58").Value = 3439.6365
$ws.Range("M58").Value = -3236.6365

$ws.Range("H99").Value = 10561.883
$ws.Range("I99").Value = 4842.143
$ws.Range("J99").Value = 12044.777
$ws.Range("K99").Value = 4842.143
$ws.Range("L99").Value = 12044.777
$ws.Range("M99").Value = -3344.143
$ws.Range("N99").Value = -15040.777

$ws.Range("H107").Value = 1853.1666
$ws.Range("J107").Value = 1750.4
$ws.Range("L107").Value = 1750.4
$ws.Range("N107").Value = -5590.4

$ws.Range("H126").Value = 10561.883
$ws.Range("I126").Value = 4842.143
$ws.Range("J126").Value = 12044.777
$ws.Range("K126").Value = 14526.429
$ws.Range("L126").Value = 36134.331
$ws.Range("M126").Value = -12056.429
$ws.Range("N126").Value = -41074.331

$ws.Range("H134").Value = 6517.037
$ws.Range("I134").Value = 4336.8237
$ws.Range("K134").Value = 13010.4711
$ws.Range("M134").Value = -10475.4711

$ws.Range("H136").Value = 5071.12
$ws.Range("I136").Value = 3439.6365
$ws.Range("K136").Value = 10318.9095
$ws.Range("M136").Value = -7768.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4916
$ws.Range("I17").Value = 5500
$ws.Range("J17").Value = 4332
$ws.Range("K17").Value = 16500
$ws.Range("L17").Value = 12996
$ws.Range("M17").Value = -16331
$ws.Range("N17").Value = -13334

$ws.Range("H120").Value = 1000
$ws.Range("I120").Value = 1000
$ws.Range("K120").Value = 3000
$ws.Range("M120").Value = 1838

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H80").Value = 2574.182
$ws.Range("I80").Value = 2314.125
$ws.Range("J80").Value = 3267.6667
$ws.Range("K80").Value = 2314.125
$ws.Range("L80").Value = 3267.6667
$ws.Range("M80").Value = -1316.125
$ws.Range("N80").Value = -5263.6667

$ws.Range("H83").Value = 2574.182
$ws.Range("I83").Value = 2314.125
$ws.Range("J83").Value = 3267.6667
$ws.Range("K83").Value = 11570.625
$ws.Range("L83").Value = 16338.3335
$ws.Range("M83").Value = -6578.625
$ws.Range("N83").Value = -26322.3335

$ws.Range("H107").Value = 534.17645
$ws.Range("J107").Value = 492.8
$ws.Range("L107").Value = 492.8
$ws.Range("N107").Value = -4332.8

$ws.Range("H113").Value = 102175.2
$ws.Range("I113").Value = 113311.445
$ws.Range("J113").Value = 1949
$ws.Range("K113").Value = 113311.445
$ws.Range("L113").Value = 1949
$ws.Range("M113").Value = -111141.445
$ws.Range("N113").Value = -6289

$ws.Range("H122").Value = 3800.3333
$ws.Range("I122").Value = 4088
$ws.Range("K122").Value = 12264
$ws.Range("M122").Value = -9814

$ws.Range("H126").Value = 3997.8696
$ws.Range("I126").Value = 3172.5
$ws.Range("K126").Value = 9517.5
$ws.Range("M126").Value = -7047.5

$ws.Range("H132").Value = 14944.875
$ws.Range("I132").Value = 10263.5
$ws.Range("K132").Value = 30790.5
$ws.Range("M132").Value = -28260.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10733.267
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 13666.333
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 13666.333
$ws.Range("M2").Value = -9888
$ws.Range("N2").Value = -13890.333

$ws.Range("H46").Value = 1700.4667
$ws.Range("I46").Value = 1037.2727
$ws.Range("J46").Value = 2084.4211
$ws.Range("K46").Value = 1037.2727
$ws.Range("L46").Value = 2084.4211
$ws.Range("M46").Value = -849.2727
$ws.Range("N46").Value = -2460.4211

$ws.Range("H61").Value = 3250.6
$ws.Range("I61").Value = 2396.3333
$ws.Range("J61").Value = 4532
$ws.Range("K61").Value = 2396.3333
$ws.Range("L61").Value = 4532
$ws.Range("M61").Value = -2194.3333
$ws.Range("N61").Value = -4936

$ws.Range("H68").Value = 3583
$ws.Range("I68").Value = 2999.8
$ws.Range("J68").Value = 6499
$ws.Range("K68").Value = 2999.8
$ws.Range("L68").Value = 6499
$ws.Range("M68").Value = -2250.8
$ws.Range("N68").Value = -7997

$ws.Range("H71").Value = 3583
$ws.Range("I71").Value = 2999.8
$ws.Range("J71").Value = 6499
$ws.Range("K71").Value = 14999
$ws.Range("L71").Value = 32495
$ws.Range("M71").Value = -11255
$ws.Range("N71").Value = -39983

$ws.Range("H82").Value = 1317.9697
$ws.Range("I82").Value = 1321.5333
$ws.Range("J82").Value = 1315
$ws.Range("K82").Value = 1321.5333
$ws.Range("L82").Value = 1315
$ws.Range("M82").Value = -960.5333000000001
$ws.Range("N82").Value = -2037

$ws.Range("H85").Value = 1317.9697
$ws.Range("I85").Value = 1321.5333
$ws.Range("J85").Value = 1315
$ws.Range("K85").Value = 1321.5333
$ws.Range("L85").Value = 1315
$ws.Range("M85").Value = -73.53330000000005
$ws.Range("N85").Value = -3811

$ws.Range("H100").Value = 5533.8
$ws.Range("I100").Value = 3499.5
$ws.Range("K100").Value = 3499.5
$ws.Range("M100").Value = -2958.5

$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490

$ws.Range("H113").Value = 3250.6
$ws.Range("I113").Value = 2396.3333
$ws.Range("J113").Value = 4532
$ws.Range("K113").Value = 2396.3333
$ws.Range("L113").Value = 4532
$ws.Range("M113").Value = -226.3332999999998
$ws.Range("N113").Value = -8872

$ws.Range("H122").Value = 4389.7334
$ws.Range("I122").Value = 5833.1665
$ws.Range("J122").Value = 3427.4443
$ws.Range("K122").Value = 17499.4995
$ws.Range("L122").Value = 10282.3329
$ws.Range("M122").Value = -15049.4995
$ws.Range("N122").Value = -15182.3329

$ws.Range("H132").Value = 4666.2925
$ws.Range("I132").Value = 4697.933
$ws.Range("K132").Value = 14093.799
$ws.Range("M132").Value = -11563.799

$ws.Range("H136").Value = 3146.3447
$ws.Range("J136").Value = 4870.2144
$ws.Range("L136").Value = 14610.6432
$ws.Range("N136").Value = -19710.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H49").Value = 99999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 99999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 99999
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -100459

$ws.Range("H75").Value = 59999
$ws.Range("I75").Value = 59999
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 59999
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -59063
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 59999
$ws.Range("I78").Value = 59999
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 179997
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -175317
$ws.Range("N78").ClearContents()

$ws.Range("H107").Value = 5849860.5
$ws.Range("I107").Value = 1257.1538
$ws.Range("J107").Value = 18521834
$ws.Range("K107").Value = 3771.4614
$ws.Range("L107").Value = 55565502
$ws.Range("M107").Value = -1851.4614
$ws.Range("N107").Value = -55569342

$ws.Range("H113").Value = 1118.2632
$ws.Range("I113").Value = 353.73334
$ws.Range("J113").Value = 3985.25
$ws.Range("K113").Value = 1061.20002
$ws.Range("L113").Value = 11955.75
$ws.Range("M113").Value = 1108.79998
$ws.Range("N113").Value = -16295.75

$ws.Range("H122").Value = 7494.0625
$ws.Range("I122").Value = 5076.4165
$ws.Range("J122").Value = 14747
$ws.Range("K122").Value = 15229.2495
$ws.Range("L122").Value = 44241
$ws.Range("M122").Value = -12779.2495
$ws.Range("N122").Value = -49141

$ws.Range("H126").Value = 9374.861999999999
$ws.Range("J126").Value = 21032
$ws.Range("L126").Value = 63096
$ws.Range("N126").Value = -68036

$ws.Range("H132").Value = 167261.47
$ws.Range("I132").Value = 274680.38
$ws.Range("K132").Value = 824041.14
$ws.Range("M132").Value = -821511.14

$ws.Range("H136").Value = 5407884.5
$ws.Range("I136").Value = 8697481
$ws.Range("J136").Value = 3547.7144
$ws.Range("K136").Value = 26092443
$ws.Range("L136").Value = 10643.1432
$ws.Range("M136").Value = -26089893
$ws.Range("N136").Value = -15743.1432
